$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resistor R7 value corrected: 560K -> 576K (matches MFR-25FBF52-576K part already in F7) ---
$ws.Range("B7").Value = "576K"

# --- TO220BH package rows: Parts reference shift as T1 (TRIAC) is renumbered into the Q-series ---
$ws.Range("E15").Value = "Q1"
$ws.Range("E16").Value = "Q2"
$ws.Range("E17").Value = "Q3"

# --- Fuse note rewritten ---
$ws.Range("H22").Value = "Derated 25% per standard, consider fast-blow?"

# --- New row 24: board heatsink / tactile switches ---
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "SWITCH-MOMENTARY-6MM"
$ws.Range("C24").Value = "SWITCH-MOMENTARY-6MM"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("D24").Value = "TACTILE-SWITCH-6MM"
$ws.Range("E24").Value = "S1, S2, S3, S4"
$ws.Range("F24").Value = "B3F-1020"

# --- Version note in row 26 ---
$ws.Range("B26").Value = "version 3.0.1"

# --- Column width tweaks (Value/Device columns widened, Parts column widened) ---
$ws.Columns.Item(2).ColumnWidth = 25.1666666666667
$ws.Columns.Item(3).ColumnWidth = 25.1666666666667
$ws.Columns.Item(5).ColumnWidth = 11.1666666666667

# --- Selection moved to the newly added version note cell ---
$ws.Range("B26").Select()
